$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'der Keksessende', 'the.ART.DEF.M.SG.NOM biscuit-end.M.SG.NOM')
    ,@(3, 'der Hase mit dem Keks', 'the.ART.DEF.M.SG.NOM hare.M.SG.NOM with the.ART.DEF.M.SG.DAT biscuit.M.SG.DAT')
    ,@(4, 'der hase ganz vorne', 'the.ART.DEF.M.SG.NOM hare.M.SG.NOM quite in-front')
    ,@(5, 'der vordere', 'the.ART.DEF.M.SG.NOM front.M.SG.NOM')
    ,@(6, 'der hase', 'the.ART.DEF.M.SG.NOM hare.M.SG.NOM')
    ,@(7, 'dass der hase blau wäre', 'that the.ART.DEF.M.SG.NOM hare.M.SG.NOM blue be.SG.Past.Sub')
    ,@(8, 'der hase', 'the.ART.DEF.M.SG.NOM hare.M.SG.NOM')
    ,@(9, 'der hase der hinfällt', 'the.ART.DEF.M.SG.NOM hare.M.SG.NOM the.DEM.M.SG.NOM fall.SG.PRS.IND')
    ,@(10, 'der Ball', 'the.ART.DEF.M.SG.NOM ball.M.SG.NOM')
    ,@(11, 'rollt', 'roll')
    ,@(12, 'der Ball, der rollt', 'the.ART.DEF.M.SG.NOM ball.M.SG.NOM  the.REL.M.SG.NOM roll.SG.PRS.IND')
    ,@(13, 'das Pferd mit dem Zopf', 'the.ART.DEF.N.SG.ACC horse.N.SG.NOM with the.ART.DEF.M.SG.DAT plait.M.SG.DAT')
    ,@(14, 'das Pferd, das die Kühe zieht', 'the.ART.DEF.N.SG.ACC horse.N.SG.ACC  the.REL.N.SG.NOM the.ART.DEF.F.PL.ACC cows.F.PL.ACC pull.SG.PRS.IND')
    ,@(15, 'da sitzen', 'there sit')
    ,@(16, 'Ziege', 'goat.SG.NOM')
    ,@(17, 'Baum', 'tree.M.SG.NOM')
    ,@(18, 'da, wo die Ziege sitzt', 'there  where.INT the.ART.DEF.F.SG.NOM goat.F.SG.NOM sit.SG.PRS.IND')
    ,@(19, 'die Frösche, die die Ente fangen', 'the.ART.DEF.M.PL.NOM frog.M.PL.NOM  the.REL.M.PL.NOM the.ART.DEF.F.SG.ACC duck.F.SG.ACC catch.PL.PRS.IND')
    ,@(20, 'die Hunde, die das Schaf gefangen genommen hat', 'the.ART.DEF.M.PL.ACC dog.M.PL.ACC  the.REL.M.PL.NOM the.ART.DEF.N.SG.ACC sheep.N.SG.ACC caught take have.SG.PRS.IND')
    ,@(21, 'der Hund, wo das Huhn den Ball zum Hund wirft', 'the.ART.DEF.M.SG.NOM dog.M.SG.NOM  where.INT the.ART.DEF.N.SG.NOM chicken.N.SG.NOM the.ART.DEF.M.SG.ACC ball.M.SG.ACC to.M.SG.DAT dog.M.SG.DAT throw.SG.PRS.IND')
    ,@(22, 'der Hund, der die Schafe wäscht', 'the.ART.DEF.M.SG.NOM dog.M.SG.NOM  the.REL.M.SG.NOM the.ART.DEF.N.PL.ACC sheep.N.PL.ACC wash.SG.PRS.IND')
    ,@(23, 'das Huhn, dem der Apfel zugeworfen wird', 'the.ART.DEF.N.SG.NOM chicken.N.SG.ACC  the.REL.M.SG.DAT the.ART.DEF.M.SG.NOM apple.M.SG.NOM thrown-at become.SG.PRS.IND')
    ,@(24, 'der Löwe, der die zwei Kühe zieht', 'the.ART.DEF.M.SG.NOM lion.M.SG.NOM  the.REL.M.SG.NOM the.ART.DEF.F.PL.ACC two cows.F.PL.ACC pull.SG.PRS.IND')
    ,@(25, 'das Huhn, das neben der Maus steht, und die Maus gibt dem Huhn die Pflaume', 'the.ART.DEF.N.SG.NOM chicken.N.SG.NOM  the.REL.N.SG.NOM next-to the.ART.DEF.F.SG.DAT mouse.F.SG.DAT stand.SG.PRS.IND  and the.ART.DEF.F.SG.NOM mouse.F.SG.NOM give.SG.PRS.IND the.ART.DEF.M.SG.DAT chicken.M.SG.DAT the.ART.DEF.F.SG.ACC plum.F.SG.ACC')
    ,@(26, 'die Pferde, die das Schaf schubst ja die und ich vermute, die Pferde, die vom Stier oder Schaf geschubst werden', 'the.ART.DEF.N.PL.ACC horse.N.PL.ACC  the.REL.N.PL.NOM the.ART.DEF.N.SG.ACC sheep.N.SG.ACC push.SG.PRS.IND yes the.DEM.F.SG.ACC and i.Prs.SG.NOM assume.SG.PRS.IND  the.ART.DEF.N.PL.ACC horse.N.PL.ACC  the.REL.N.PL.NOM from.M.SG.DAT bull.M.SG.DAT or sheep.N.SG.DAT push become.PL.PRS.IND')
    ,@(27, 'der Frosch, der da sitzt unter dem Tisch', 'the.ART.DEF.M.SG.NOM frog.M.SG.NOM  the.REL.M.SG.NOM there sit.SG.PRS.IND under the.ART.DEF.M.SG.DAT table.M.SG.DAT')
    ,@(28, 'der Tisch', 'the.ART.DEF.M.SG.NOM table.M.SG.NOM')
    ,@(29, 'der, wo der Frosch drunter sitzt', 'the.DEM.M.SG.NOM  where.INT the.ART.DEF.M.SG.NOM frog.M.SG.NOM underneath sit.SG.PRS.IND')
    ,@(30, 'die Ziegen, die die Katze schubsen', 'the.ART.DEF.F.PL.ACC goat.F.PL.ACC  the.REL.F.PL.NOM the.ART.DEF.F.SG.ACC cat.F.SG.ACC push.PL.PRS.IND')
    ,@(31, 'die Katze, die vor der Ente steht und die Ente gibt der Katze die Blume', 'the.ART.DEF.F.SG.NOM cat.F.SG.NOM  the.REL.F.SG.NOM before the.ART.DEF.F.SG.DAT duck.F.SG.DAT stand.SG.PRS.IND and the.ART.DEF.F.SG.NOM duck.F.SG.NOM give.SG.PRS.IND the.ART.DEF.F.SG.DAT cat.F.SG.DAT the.ART.DEF.F.SG.ACC flower.F.SG.ACC')
    ,@(32, 'die Katze', 'the.ART.DEF.F.SG.NOM cat.F.SG.NOM')
    ,@(33, 'der Frosch, der von den Enten gewaschen wird', 'the.ART.DEF.M.SG.NOM frog.M.SG.NOM  the.REL.M.SG.NOM from the.ART.DEF.F.PL.DAT duck.F.PL.DAT wash become.SG.PRS.IND')
    ,@(34, 'das Bett, auf dem die Katze springt', 'the.ART.DEF.N.SG.ACC bed.N.SG.ACC  on the.REL.N.SG.DAT the.ART.DEF.F.SG.NOM cat.F.SG.NOM jump.SG.PRS.IND')
    ,@(35, 'der Löwe, der von den Katzen gezogen wird', 'the.ART.DEF.M.SG.NOM lion.M.SG.NOM  the.REL.M.SG.NOM from the.ART.DEF.F.PL.DAT cat.F.PL.DAT pull become.SG.PRS.IND')
    ,@(36, 'das Haus, auf dem das Schaf springt', 'the.ART.DEF.N.SG.ACC house.N.SG.ACC  on the.REL.N.SG.DAT the.ART.DEF.N.SG.NOM sheep.N.SG.NOM jump.SG.PRS.IND')
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
}
